$wb = $excel.ActiveWorkbook

# Add the two new sheets in the right position
$generic = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$generic.Name = "generic_text_elements"

$case = $wb.Worksheets.Add($wb.Worksheets.Item(3))
$case.Name = "case_text_elements"

# Populate generic_text_elements
$generic.Range("A1").Value = "generic_text_element"
$generic.Range("B1").Value = "value"

$generic.Range("A2").Value = "title_strategic_challenge"
$generic.Range("B2").Value = "Strategic Challenge"

$generic.Range("A3").Value = "title_key_outputs"
$generic.Range("B3").Value = "Key outputs"

$generic.Range("A4").Value = "title_dmo"
$generic.Range("B4").Value = "Options"

$generic.Range("A5").Value = "title_scenarios"
$generic.Range("B5").Value = "Scenarios"

$generic.Range("A6").Value = "title_comparison"
$generic.Range("B6").Value = "Comparisons of options"

$generic.Range("A7").Value = "title_theme_weights"
$generic.Range("B7").Value = "Key output and theme weights"

$generic.Range("A8").Value = "title_scenario_weights"
$generic.Range("B8").Value = "Scenario weights"

$generic.Range("A9").Value = "text_strategic_challenge"
$generic.Range("B9").Value = "Describing strategic challenge that requires a decision"

$generic.Range("A10").Value = "text_key_outputs"
$generic.Range("B10").Value = "Which indicators do you use to evaluate the impact of your decision(s)?"

$generic.Range("A11").Value = "text_dmo"
$generic.Range("B11").Value = "Which options do you have to influence your impact?"

$generic.Range("A12").Value = "text_scenarios"
$generic.Range("B12").Value = "Which uncertainty do you want to account for?"

# Populate case_text_elements
$case.Range("A1").Value = "case_text_element"
$case.Range("B1").Value = "value"

$case.Range("A2").Value = "strategic_challenge"
$case.Range("B2").Value = "How to source energy?"

"done"
